$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains text formatting so numeric-looking strings
# like "26.267.07" or "1.008" are not converted to actual numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.267.07"
$ws.Range("E2").Value = "  -0.26%  "

$ws.Range("D3").Value = "1.688.80"
$ws.Range("E3").Value = "  +0.54%  "

$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "217.36"
$ws.Range("E5").Value = "  -0.31%  "

$ws.Range("D6").Value = "0.5353"
$ws.Range("E6").Value = "  +1.57%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "0.2716"
$ws.Range("E8").Value = "  +0.81%  "

$ws.Range("D9").Value = "0.06413"
$ws.Range("E9").Value = "  -0.86%  "

$ws.Range("D10").Value = "21.67"
$ws.Range("E10").Value = "  -1.32%  "

$ws.Range("D11").Value = "0.07669"
$ws.Range("E11").Value = "  +2.13%  "

$ws.Range("D12").Value = "1.694.02"
$ws.Range("E12").Value = "  +0.67%  "

$ws.Range("D13").Value = "4.512"
$ws.Range("E13").Value = "  -0.09%  "

$ws.Range("D14").Value = "0.5777"
$ws.Range("E14").Value = "  -0.21%  "

$ws.Range("D15").Value = "0.000008330"
$ws.Range("E15").Value = "  -2.21%  "

$ws.Range("D16").Value = "66.43"
$ws.Range("E16").Value = "  +2.55%  "

$ws.Range("D17").Value = "26.308.24"
$ws.Range("E17").Value = "  -0.10%  "

$ws.Range("D18").Value = "1.008"
$ws.Range("E18").Value = "  +0.05%  "

$ws.Range("D19").Value = "4.887"
$ws.Range("E19").Value = "  -0.80%  "

$ws.Range("D20").Value = "10.84"
$ws.Range("E20").Value = "  -0.38%  "

$ws.Range("D21").Value = "190.09"
$ws.Range("E21").Value = "  +0.16%  "

$ws.Range("D22").Value = "6.236"
$ws.Range("E22").Value = "  +0.54%  "

$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("D24").Value = "148.69"
$ws.Range("E24").Value = "  +2.63%  "

$ws.Range("D25").Value = "0.1283"
$ws.Range("E25").Value = "  +2.13%  "

$ws.Range("D26").Value = "7.823"
$ws.Range("E26").Value = "  +0.43%  "

$ws.Range("D27").Value = "15.74"
$ws.Range("E27").Value = "  -0.21%  "

$ws.Range("D28").Value = "1.378"
$ws.Range("E28").Value = "  +1.20%  "

$ws.Range("D29").Value = "0.06133"
$ws.Range("E29").Value = "  -5.68%  "

$ws.Range("D30").Value = "1.322"
$ws.Range("E30").Value = "  -0.15%  "

$ws.Range("D31").Value = "3.579"
$ws.Range("E31").Value = "  -0.34%  "

$ws.Range("D32").Value = "3.581"
$ws.Range("E32").Value = "  -0.19%  "

$ws.Range("D33").Value = "1.682"
$ws.Range("E33").Value = "  +1.30%  "

$ws.Range("D34").Value = "1.031"
$ws.Range("E34").Value = "  +0.19%  "

$ws.Range("D35").Value = "0.6177"
$ws.Range("E35").Value = "  -0.59%  "

$ws.Range("D36").Value = "2.429"
$ws.Range("E36").Value = "  +0.86%  "

$ws.Range("E37").Value = "  +0.91%  "

$ws.Range("E38").Value = "  +1.23%  "

$ws.Range("D39").Value = "6.158"
$ws.Range("E39").Value = "  -2.19%  "

$ws.Range("D40").Value = "1.104.92"
$ws.Range("E40").Value = "  -0.99%  "

$ws.Range("D41").Value = "0.8780"
$ws.Range("E41").Value = "  +0.45%  "

$ws.Range("E42").Value = "  -0.35%  "

$ws.Range("D43").Value = "100.79"
$ws.Range("E43").Value = "  +0.32%  "

$ws.Range("D44").Value = "1.839.11"
$ws.Range("E44").Value = "  +0.61%  "

$ws.Range("D45").Value = "0.00000000111"
$ws.Range("E45").Value = "  +0.76%  "

$ws.Range("D46").Value = "57.57"
$ws.Range("E46").Value = "  +1.13%  "

$ws.Range("D47").Value = "1.007"
$ws.Range("E47").Value = "  +1.03%  "

$ws.Range("D48").Value = "8.103"
$ws.Range("E48").Value = "  -0.60%  "

$ws.Range("D49").Value = "0.05285"
$ws.Range("E49").Value = "  +0.26%  "

$ws.Range("D50").Value = "0.4294"
$ws.Range("E50").Value = "  +0.00%  "

$ws.Range("D51").Value = "6.049"
$ws.Range("E51").Value = "  -0.49%  "

